$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 5) down onto the
# new row (row 6) so the new row gets identical cell styles (borders,
# fonts, number formats) to the rest of the table.
$ws.Range("A5:S5").Copy()
$ws.Range("A6:S6").PasteSpecial(-4122)

# --- New session / speaker row (Narek Babajanyan) ---
$ws.Range("A6").Value = 1061298
$ws.Range("B6").Value = "OAuth supply-chain risks: How to mitigate them, and how not to become them"
$ws.Range("C6").Value = "A lightning session about Google APIs' OAuth permission scopes and their categories - from a security engineer's perspective. `n Topics include an introduction to supply chain risks, OAuth illicit consent attacks, `n The talk intends to be an awareness session for developers to request the least possible amount of permissions (and data) and for security engineers to audit their Google Workspace OAuth integrations."
$ws.Range("D6").Value = "Narek Babajanyan"
$ws.Range("E6").Value = "narek_babajanyan@outlook.com"
$ws.Range("F6").Value = 45960.46875
$ws.Range("G6").Value = "No"
$ws.Range("H6").Value = "Hall A"
$ws.Range("I6").Value = 46011.5
$ws.Range("J6").Value = 10
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = "70d1cacf-57c0-4903-ab30-ef0faf6f1955"
$ws.Range("N6").Value = "Narek"
$ws.Range("O6").Value = "Babajanyan"
$ws.Range("P6").Value = "narek_babajanyan@outlook.com"
$ws.Range("Q6").Value = "Cybersecurity @ ServiceTitan | Armed Forces ex-officer"
$ws.Range("R6").Value = "I am currently an Incident Response and Threat Prevention Engineer at ServiceTitan. Previously, I helped safeguard Armenia's critical infrastructure at the Information Systems Agency of Armenia."
$ws.Range("S6").Value = "https://sessionize.com/image/71f5-400o400o1-BVLaZdEqXxPwDdXZowmvfo.png"

# Profile picture column (S) carries an external hyperlink for every row,
# so add one for the newly added row too.
$ws.Hyperlinks.Add($ws.Range("S6"), "https://sessionize.com/image/71f5-400o400o1-BVLaZdEqXxPwDdXZowmvfo.png")

# Adding the hyperlink re-styles the cell with a generic "Hyperlink" look;
# restore the table's own underlined-link style (matching S2:S5) so S6
# stays visually consistent with the rest of the column.
$ws.Range("S5").Copy()
$ws.Range("S6").PasteSpecial(-4122)
